$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.979.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.53%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.822.44'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.77%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.010'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.59%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.87%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6143'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.39%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.008'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.25%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07328'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2888'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.29%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.68%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07685'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.69%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.806.96'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.47%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.931'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.19%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6618'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.14%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '81.74'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.13%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008947'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.85%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.831'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.83%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '28.941.28'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.76%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.075.79'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.08%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '234.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.83%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.27%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.008'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.12%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.074'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.79%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.014'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.84%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.90'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.99%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1394'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.23%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.421'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.97%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.64'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.494'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.53%  '

$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05535'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.81%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.081'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.090'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.60%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.208'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.14%  '

$ws.Range("E34").Value = '  -1.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7312'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.61%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.129'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.76%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.648'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.812'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.47%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01760'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.72%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.190.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.324'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.61%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8962'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.36%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.007'
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.63'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.51%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.954.91'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.17%  '

$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5123'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.30%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '64.03'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.54%  '

$ws.Range("E48").Value = '  -1.89%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.3986'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.022'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.20%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05806'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.63%  '
